# Lit Review Matrix -> expand into tidy data:
#   matrix (renamed from Sheet1), data, policy group, outcome groups
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the original sheet and add the three new sheets, in tab order.
# ---------------------------------------------------------------------------
$matrix = $wb.Worksheets.Item(1)
$matrix.Name = "matrix"

$data = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$data.Name = "data"

$policyGroup = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$policyGroup.Name = "policy group"

$outcomeGroups = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$outcomeGroups.Name = "outcome groups"

# ---------------------------------------------------------------------------
# 2. "matrix" sheet: cursor left on F1 once the restructuring was done.
# ---------------------------------------------------------------------------
$matrix.Range("F1").Select()

# ---------------------------------------------------------------------------
# 3. "data" sheet headers (Paper / Policy / Outcome first, Results / Summary
#    filled in later once the two "group" sheets exist).
# ---------------------------------------------------------------------------
$data.Range("B2").Value = "Paper"
$data.Range("C2").Value = "Policy"
$data.Range("D2").Value = "Outcome"

# ---------------------------------------------------------------------------
# 4. "policy group" sheet: one row per policy named on "matrix", column A.
# ---------------------------------------------------------------------------
$policyGroup.Range("B2").Value = "Policy"
$policyGroup.Range("C2").Value = "Policy Group"

$policies = @(
    "Policies related to childcare (e.g., access, affordability, etc.)",
    "Employment, unemployment, and labor policies ",
    "Housing policies (e.g., mortgage access and assistance, property tax policies, eviction and displacement protection, rent protection, etc.)",
    "Income inequity and poverty",
    "Investments in public schooling",
    "Insurance access (public or private)",
    "Paid leave (e.g., family leave, medical/sick leave)",
    "Policies related to the criminal legal system, policing, incarceration, and re-entry",
    "Policies related to racial discrimination (e.g., Jim Crow, redlining, etc.)",
    "Policies related to community investment (or disinvestment) and neighborhood advantage",
    "Policies related to wages and economic inclusion (e.g., living wages, wage increases, Universal Basic Income, etc.)",
    "Policies related to zoning, land use, street design, and community design (e.g., Complete Streets, inclusionary zoning, etc.)",
    "Tax policy (e.g., how the local and federal tax system funds schools, influences wealth, impacts communities of color, etc.)",
    "Voting access, voting rights, and civic participation policies",
    "Others?"
)
for ($i = 0; $i -lt $policies.Length; $i++) {
    $row = 3 + $i
    $policyGroup.Cells.Item($row, 2).Value = $policies[$i]
}

# ---------------------------------------------------------------------------
# 5. "outcome groups" sheet: one row per outcome named on "matrix", row 1.
# ---------------------------------------------------------------------------
$outcomeGroups.Range("B2").Value = "Outcome"
$outcomeGroups.Range("C2").Value = "Outcome Group"

$outcomes = @(
    "Chronic diseases",
    "Homicide",
    "Infant Mortality",
    "Injuries (Unintentional)",
    "Life Expectancy / Premature Mortality",
    "Mental health conditions",
    "Poisoning / Overdoses",
    "Sexually transmitted infections",
    "Substance use disorders",
    "Suicide"
)
for ($i = 0; $i -lt $outcomes.Length; $i++) {
    $row = 3 + $i
    $outcomeGroups.Cells.Item($row, 2).Value = $outcomes[$i]
}

# ---------------------------------------------------------------------------
# 6. Back to "data": finish the headers, then the 3 literature rows.
# ---------------------------------------------------------------------------
$jimCrowResults = "Jim Crow birthplace was associated with increased odds of ER" + [char]0x2212 + " breast cancer only among the black, not white women, with the effect strongest for women born before 1965. Among black women, the odds ratio (OR) for an ER" + [char]0x2212 + " tumor, comparing women born in a Jim Crow versus not Jim Crow state, equaled 1.09 (95% confidence interval [CI] 1.06, 1.13), on par with the OR comparing women in the worst versus best census tract socioeconomic quintiles (1.15; 95% CI 1.07, 1.23). The black versus white OR for ER" + [char]0x2212 + " was higher among women born in Jim Crow versus non-Jim Crow states (1.41 [95% CI 1.13, 1.46] vs. 1.27 [95% CI 1.24, 1.31])."

$data.Range("E2").Value = "Results"
$data.Range("E3").Value = $jimCrowResults
$data.Range("F2").Value = "Summary"

$data.Range("B3").Value = "Krieger et al. (2017): https://link.springer.com/article/10.1007%2Fs10552-016-0834-2"
$data.Range("C3").Value = "Policies related to racial discrimination (e.g., Jim Crow, redlining, etc.)"
$data.Range("D3").Value = "Chronic diseases"

$data.Range("B4").Value = "Krieger et al. (2014): https://www.ncbi.nlm.nih.gov/pmc/articles/PMC3828968/"
$data.Range("C4").Value = "Policies related to racial discrimination (e.g., Jim Crow, redlining, etc.)"
$data.Range("D4").Value = "Infant Mortality"

$data.Range("B5").Value = "Krieger et al. (2014): https://www.ncbi.nlm.nih.gov/pmc/articles/PMC4710482/"
$data.Range("C5").Value = "Policies related to racial discrimination (e.g., Jim Crow, redlining, etc.)"
$data.Range("D5").Value = "Life Expectancy / Premature Mortality"

# ---------------------------------------------------------------------------
# 7. Formatting.
# ---------------------------------------------------------------------------

# -- data --
$data.Range("B2:F2").Font.Name = "Calibri"
$data.Range("B2:F2").Font.Size = 11
$data.Range("B2:F2").Font.Bold = $true
$data.Range("B2:F2").HorizontalAlignment = -4131
$data.Range("B2:F2").VerticalAlignment = -4160

$data.Range("B3:F6").HorizontalAlignment = -4131
$data.Range("B3:F6").VerticalAlignment = -4160
$data.Range("B3:E5").WrapText = $true

$data.Range("B2:F6").Borders.LineStyle = 1

$data.Rows.Item(3).RowHeight = 120
$data.Rows.Item(4).RowHeight = 45
$data.Rows.Item(5).RowHeight = 45

$data.Columns.Item(2).ColumnWidth = 40.6
$data.Columns.Item(3).ColumnWidth = 36.45
$data.Columns.Item(4).ColumnWidth = 19.45
$data.Columns.Item(5).ColumnWidth = 81.15
$data.Columns.Item(6).ColumnWidth = 21.73

# -- policy group --
$policyGroup.Range("B2:C2").Font.Name = "Calibri"
$policyGroup.Range("B2:C2").Font.Size = 11
$policyGroup.Range("B2:C2").Font.Bold = $true

$policyGroup.Range("B3:B17").Font.Name = "Arial"
$policyGroup.Range("B3:B17").Font.Size = 10
$policyGroup.Range("B3:B17").Font.Bold = $true
$policyGroup.Range("B3:B17").WrapText = $true

$policyGroup.Range("B2:C17").Borders.LineStyle = 1

$policyGroup.Range("B5").RowHeight = 39
$policyGroup.Rows.Item(5).RowHeight = 39
$policyGroup.Rows.Item(11).RowHeight = 26.25
$policyGroup.Rows.Item(3).RowHeight = 26.25
$policyGroup.Rows.Item(10).RowHeight = 26.25
$policyGroup.Rows.Item(12).RowHeight = 26.25
$policyGroup.Rows.Item(13).RowHeight = 39
$policyGroup.Rows.Item(14).RowHeight = 39
$policyGroup.Rows.Item(15).RowHeight = 39
$policyGroup.Rows.Item(16).RowHeight = 26.25

$policyGroup.Columns.Item(2).ColumnWidth = 44
$policyGroup.Columns.Item(3).ColumnWidth = 18.29

# -- outcome groups --
$outcomeGroups.Range("B2").Font.Name = "Calibri"
$outcomeGroups.Range("B2").Font.Size = 11
$outcomeGroups.Range("B2").Font.Bold = $true
$outcomeGroups.Range("B2").HorizontalAlignment = -4108

$outcomeGroups.Range("C2").Font.Name = "Calibri"
$outcomeGroups.Range("C2").Font.Size = 11
$outcomeGroups.Range("C2").Font.Bold = $true

$outcomeGroups.Range("B3:B12").HorizontalAlignment = -4131
$outcomeGroups.Range("B3:B12").VerticalAlignment = -4160

$outcomeGroups.Range("B2:C12").Borders.LineStyle = 1

$outcomeGroups.Range("B6").RowHeight = 25.5
$outcomeGroups.Rows.Item(6).RowHeight = 25.5
$outcomeGroups.Rows.Item(7).RowHeight = 25.5
$outcomeGroups.Rows.Item(8).RowHeight = 25.5
$outcomeGroups.Rows.Item(9).RowHeight = 25.5
$outcomeGroups.Rows.Item(10).RowHeight = 25.5

$outcomeGroups.Columns.Item(2).ColumnWidth = 19.43
$outcomeGroups.Columns.Item(3).ColumnWidth = 19.43

# ---------------------------------------------------------------------------
# 8. Selections + active sheet (last edited == "outcome groups").
# ---------------------------------------------------------------------------
$data.Range("F3").Select()
$policyGroup.Range("B5:B8").Select()
$outcomeGroups.Range("B3").Select()
$outcomeGroups.Activate()
